$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for column P (year 2022), keyed by row number
$values = @{
    4  = 2022
    5  = 11.4
    6  = 12.6
    7  = 9.8
    8  = 11.4
    9  = 5.4
    10 = 4.7
    11 = 3.4
    12 = 17.7
    13 = 20.5
    14 = 8.4
    16 = 12.9
    17 = 10.5
}

foreach ($row in $values.Keys) {
    $src = $ws.Range("O" + $row)
    $dst = $ws.Range("P" + $row)
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = $values[$row]
}

$excel.CutCopyMode = 0

# Update selected cell as recorded in the saved view
$ws.Range("Q4").Select()
